$wb = $excel.ActiveWorkbook

# --- Rename existing sheets' MATCH_CARD_LINK column to MATCH_CODE and
#     replace the full scorecard URL with the bare match code -----------

$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4485"

$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4485"

# --- Insert a new "Player Info" sheet before "ODI Batting" -------------

$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4668"
$playerInfo.Range("B2").Value = "Krishnappa Gowtham"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
